$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.514509333333333
$ws.Range("H2").Value = 4.543528
$ws.Range("I2").Value = 0.01996786707219448
$ws.Range("J2").Value = 0.02165084619119693
$ws.Range("M2").Value = 28.85518433333334
$ws.Range("N2").Value = 86.56555300000001
$ws.Range("O2").Value = 0.1999651185353207
$ws.Range("P2").Value = 0.2044513327926365
$ws.Range("Q2").Value = 43.70144598788712
$ws.Range("R2").Value = 393.3130138909841
$ws.Range("S2").Value = 0.003992876905988895
$ws.Range("T2").Value = 0.004426544359878591
$ws.Range("G3").Value = 1.514509333333333
$ws.Range("H3").Value = 4.543528
$ws.Range("I3").Value = 0.01996786707219448
$ws.Range("J3").Value = 0.02165084619119693
$ws.Range("O3").Value = 0.3546352265743414
$ws.Range("P3").Value = 0.3625914622481308
$ws.Range("Q3").Value = 77.50387824166046
$ws.Range("R3").Value = 697.534904174944
$ws.Range("S3").Value = 0.00708130906335402
$ws.Range("T3").Value = 0.007850411979375471
$ws.Range("G4").Value = 1.514509333333333
$ws.Range("H4").Value = 4.543528
$ws.Range("I4").Value = 0.01996786707219448
$ws.Range("J4").Value = 0.02165084619119693
$ws.Range("M4").Value = 29.393479
$ws.Range("N4").Value = 88.180437
$ws.Range("O4").Value = 0.2036954761578358
$ws.Range("P4").Value = 0.2082653809291453
$ws.Range("Q4").Value = 44.51669828463734
$ws.Range("R4").Value = 400.650284561736
$ws.Range("S4").Value = 0.004067364191127025
$ws.Range("T4").Value = 0.004509121729447963
$ws.Range("G5").Value = 1.514509333333333
$ws.Range("H5").Value = 4.543528
$ws.Range("I5").Value = 0.01996786707219448
$ws.Range("J5").Value = 0.02165084619119693
$ws.Range("M5").Value = 9.499066500000001
$ws.Range("N5").Value = 18.998133
$ws.Range("O5").Value = 0.0658280999596015
$ws.Range("P5").Value = 0.04486996822421697
$ws.Range("Q5").Value = 14.386424872204
$ws.Range("R5").Value = 86.31854923322402
$ws.Range("S5").Value = 0.001314446749608453
$ws.Range("T5").Value = 0.0009714727806264154
$ws.Range("G6").Value = 1.514509333333333
$ws.Range("H6").Value = 4.543528
$ws.Range("I6").Value = 0.01996786707219448
$ws.Range("J6").Value = 0.02165084619119693
$ws.Range("M6").Value = 25.37910966666666
$ws.Range("N6").Value = 76.13732899999999
$ws.Range("O6").Value = 0.1758760787729007
$ws.Range("P6").Value = 0.1798218558058706
$ws.Range("Q6").Value = 38.43689846185689
$ws.Range("R6").Value = 345.932086156712
$ws.Range("S6").Value = 0.003511870162116085
$ws.Range("T6").Value = 0.003893295341868499
$ws.Range("I7").Value = 0.6527104067845205
$ws.Range("J7").Value = 0.7077236929508544
$ws.Range("M7").Value = 28.85518433333334
$ws.Range("N7").Value = 86.56555300000001
$ws.Range("O7").Value = 0.1999651185353207
$ws.Range("P7").Value = 0.2044513327926365
$ws.Range("Q7").Value = 1428.514547131884
$ws.Range("R7").Value = 12856.63092418696
$ws.Range("S7").Value = 0.130519313861904
$ws.Range("T7").Value = 0.1446950522727288
$ws.Range("I8").Value = 0.6527104067845205
$ws.Range("J8").Value = 0.7077236929508544
$ws.Range("O8").Value = 0.3546352265743414
$ws.Range("P8").Value = 0.3625914622481308
$ws.Range("Q8").Value = 2533.449752624607
$ws.Range("S8").Value = 0.231474102997459
$ws.Range("T8").Value = 0.2566145686946975
$ws.Range("I9").Value = 0.6527104067845205
$ws.Range("J9").Value = 0.7077236929508544
$ws.Range("M9").Value = 29.393479
$ws.Range("N9").Value = 88.180437
$ws.Range("O9").Value = 0.2036954761578358
$ws.Range("P9").Value = 0.2082653809291453
$ws.Range("Q9").Value = 1455.163545561207
$ws.Range("R9").Value = 13096.47191005087
$ws.Range("S9").Value = 0.1329541571031476
$ws.Range("T9").Value = 0.1473943445049911
$ws.Range("I10").Value = 0.6527104067845205
$ws.Range("J10").Value = 0.7077236929508544
$ws.Range("M10").Value = 9.499066500000001
$ws.Range("N10").Value = 18.998133
$ws.Range("O10").Value = 0.0658280999596015
$ws.Range("P10").Value = 0.04486996822421697
$ws.Range("Q10").Value = 470.2640094988991
$ws.Range("R10").Value = 2821.584056993395
$ws.Range("S10").Value = 0.04296668590248357
$ws.Range("T10").Value = 0.03175553961423032
$ws.Range("I11").Value = 0.6527104067845205
$ws.Range("J11").Value = 0.7077236929508544
$ws.Range("M11").Value = 25.37910966666666
$ws.Range("N11").Value = 76.13732899999999
$ws.Range("O11").Value = 0.1758760787729007
$ws.Range("P11").Value = 0.1798218558058706
$ws.Range("Q11").Value = 1256.42681513588
$ws.Range("R11").Value = 11307.84133622292
$ws.Range("S11").Value = 0.1147961469195264
$ws.Range("T11").Value = 0.1272641878642068
$ws.Range("G12").Value = 3.794695333333333
$ws.Range("H12").Value = 11.384086
$ws.Range("I12").Value = 0.05003070653167101
$ws.Range("J12").Value = 0.05424751316892035
$ws.Range("M12").Value = 28.85518433333334
$ws.Range("N12").Value = 86.56555300000001
$ws.Range("O12").Value = 0.1999651185353207
$ws.Range("P12").Value = 0.2044513327926365
$ws.Range("Q12").Value = 109.4966333321731
$ws.Range("R12").Value = 985.4696999895581
$ws.Range("S12").Value = 0.01000439616201144
$ws.Range("T12").Value = 0.01109097636807187
$ws.Range("G13").Value = 3.794695333333333
$ws.Range("H13").Value = 11.384086
$ws.Range("I13").Value = 0.05003070653167101
$ws.Range("J13").Value = 0.05424751316892035
$ws.Range("O13").Value = 0.3546352265743414
$ws.Range("P13").Value = 0.3625914622481308
$ws.Range("Q13").Value = 194.1906851320365
$ws.Range("R13").Value = 1747.716166188328
$ws.Range("S13").Value = 0.01774265094653353
$ws.Range("T13").Value = 0.01966968512324356
$ws.Range("G14").Value = 3.794695333333333
$ws.Range("H14").Value = 11.384086
$ws.Range("I14").Value = 0.05003070653167101
$ws.Range("J14").Value = 0.05424751316892035
$ws.Range("M14").Value = 29.393479
$ws.Range("N14").Value = 88.180437
$ws.Range("O14").Value = 0.2036954761578358
$ws.Range("P14").Value = 0.2082653809291453
$ws.Range("Q14").Value = 111.5392975917313
$ws.Range("R14").Value = 1003.853678325582
$ws.Range("S14").Value = 0.01019102858948167
$ws.Range("T14").Value = 0.01129787899458402
$ws.Range("G15").Value = 3.794695333333333
$ws.Range("H15").Value = 11.384086
$ws.Range("I15").Value = 0.05003070653167101
$ws.Range("J15").Value = 0.05424751316892035
$ws.Range("M15").Value = 9.499066500000001
$ws.Range("N15").Value = 18.998133
$ws.Range("O15").Value = 0.0658280999596015
$ws.Range("P15").Value = 0.04486996822421697
$ws.Range("Q15").Value = 36.046063318573
$ws.Range("R15").Value = 216.276379911438
$ws.Range("S15").Value = 0.003293426350616327
$ws.Range("T15").Value = 0.002434084192132248
$ws.Range("G16").Value = 3.794695333333333
$ws.Range("H16").Value = 11.384086
$ws.Range("I16").Value = 0.05003070653167101
$ws.Range("J16").Value = 0.05424751316892035
$ws.Range("M16").Value = 25.37910966666666
$ws.Range("N16").Value = 76.13732899999999
$ws.Range("O16").Value = 0.1758760787729007
$ws.Range("P16").Value = 0.1798218558058706
$ws.Range("Q16").Value = 96.30598901625488
$ws.Range("R16").Value = 866.753901146294
$ws.Range("S16").Value = 0.008799204483028044
$ws.Range("T16").Value = 0.009754888490888663
$ws.Range("G17").Value = 17.6874565
$ws.Range("H17").Value = 35.374913
$ws.Range("I17").Value = 0.2331981536620147
$ws.Range("J17").Value = 0.1685687422615142
$ws.Range("M17").Value = 28.85518433333334
$ws.Range("N17").Value = 86.56555300000001
$ws.Range("O17").Value = 0.1999651185353207
$ws.Range("P17").Value = 0.2044513327926365
$ws.Range("Q17").Value = 510.3748176953149
$ws.Range("R17").Value = 3062.248906171889
$ws.Range("S17").Value = 0.0466314964392427
$ws.Range("T17").Value = 0.034464104022545
$ws.Range("G18").Value = 17.6874565
$ws.Range("H18").Value = 35.374913
$ws.Range("I18").Value = 0.2331981536620147
$ws.Range("J18").Value = 0.1685687422615142
$ws.Range("O18").Value = 0.3546352265743414
$ws.Range("P18").Value = 0.3625914622481308
$ws.Range("Q18").Value = 905.1423090034874
$ws.Range("R18").Value = 5430.853854020924
$ws.Range("S18").Value = 0.08270028006064668
$ws.Range("T18").Value = 0.06112158674593071
$ws.Range("G19").Value = 17.6874565
$ws.Range("H19").Value = 35.374913
$ws.Range("I19").Value = 0.2331981536620147
$ws.Range("J19").Value = 0.1685687422615142
$ws.Range("M19").Value = 29.393479
$ws.Range("N19").Value = 88.180437
$ws.Range("O19").Value = 0.2036954761578358
$ws.Range("P19").Value = 0.2082653809291453
$ws.Range("Q19").Value = 519.8958811961635
$ws.Range("R19").Value = 3119.375287176981
$ws.Range("S19").Value = 0.04750140894931226
$ws.Range("T19").Value = 0.03510703331984116
$ws.Range("G20").Value = 17.6874565
$ws.Range("H20").Value = 35.374913
$ws.Range("I20").Value = 0.2331981536620147
$ws.Range("J20").Value = 0.1685687422615142
$ws.Range("M20").Value = 9.499066500000001
$ws.Range("N20").Value = 18.998133
$ws.Range("O20").Value = 0.0658280999596015
$ws.Range("P20").Value = 0.04486996822421697
$ws.Range("Q20").Value = 168.0143255093573
$ws.Range("R20").Value = 672.0573020374291
$ws.Range("S20").Value = 0.01535099136965761
$ws.Range("T20").Value = 0.00756367410887036
$ws.Range("G21").Value = 17.6874565
$ws.Range("H21").Value = 35.374913
$ws.Range("I21").Value = 0.2331981536620147
$ws.Range("J21").Value = 0.1685687422615142
$ws.Range("M21").Value = 25.37910966666666
$ws.Range("N21").Value = 76.13732899999999
$ws.Range("O21").Value = 0.1758760787729007
$ws.Range("P21").Value = 0.1798218558058706
$ws.Range("Q21").Value = 448.8918982378961
$ws.Range("R21").Value = 2693.351389427377
$ws.Range("S21").Value = 0.04101397684315549
$ws.Range("T21").Value = 0.03031234406432697
$ws.Range("G22").Value = 3.344326
$ws.Range("H22").Value = 10.032978
$ws.Range("I22").Value = 0.04409286594959943
$ws.Range("J22").Value = 0.04780920542751418
$ws.Range("M22").Value = 28.85518433333334
$ws.Range("N22").Value = 86.56555300000001
$ws.Range("O22").Value = 0.1999651185353207
$ws.Range("P22").Value = 0.2044513327926365
$ws.Range("Q22").Value = 96.50114320075934
$ws.Range("R22").Value = 868.5102888068341
$ws.Range("S22").Value = 0.008817035166173654
$ws.Range("T22").Value = 0.009774655769412223
$ws.Range("G23").Value = 3.344326
$ws.Range("H23").Value = 10.032978
$ws.Range("I23").Value = 0.04409286594959943
$ws.Range("J23").Value = 0.04780920542751418
$ws.Range("O23").Value = 0.3546352265743414
$ws.Range("P23").Value = 0.3625914622481308
$ws.Range("Q23").Value = 171.1433725759493
$ws.Range("R23").Value = 1540.290353183544
$ws.Range("S23").Value = 0.01563688350634826
$ws.Range("T23").Value = 0.01733520970488364
$ws.Range("G24").Value = 3.344326
$ws.Range("H24").Value = 10.032978
$ws.Range("I24").Value = 0.04409286594959943
$ws.Range("J24").Value = 0.04780920542751418
$ws.Range("M24").Value = 29.393479
$ws.Range("N24").Value = 88.180437
$ws.Range("O24").Value = 0.2036954761578358
$ws.Range("P24").Value = 0.2082653809291453
$ws.Range("Q24").Value = 98.301376050154
$ws.Range("R24").Value = 884.712384451386
$ws.Range("S24").Value = 0.008981517324767282
$ws.Range("T24").Value = 0.009957002380280999
$ws.Range("G25").Value = 3.344326
$ws.Range("H25").Value = 10.032978
$ws.Range("I25").Value = 0.04409286594959943
$ws.Range("J25").Value = 0.04780920542751418
$ws.Range("M25").Value = 9.499066500000001
$ws.Range("N25").Value = 18.998133
$ws.Range("O25").Value = 0.0658280999596015
$ws.Range("P25").Value = 0.04486996822421697
$ws.Range("Q25").Value = 31.76797507167901
$ws.Range("R25").Value = 190.607850430074
$ws.Range("S25").Value = 0.002902549587235541
$ws.Range("T25").Value = 0.002145197528357623
$ws.Range("G26").Value = 3.344326
$ws.Range("H26").Value = 10.032978
$ws.Range("I26").Value = 0.04409286594959943
$ws.Range("J26").Value = 0.04780920542751418
$ws.Range("M26").Value = 25.37910966666666
$ws.Range("N26").Value = 76.13732899999999
$ws.Range("O26").Value = 0.1758760787729007
$ws.Range("P26").Value = 0.1798218558058706
$ws.Range("Q26").Value = 84.87601631508467
$ws.Range("R26").Value = 763.884146835762
$ws.Range("S26").Value = 0.007754880365074698
$ws.Range("T26").Value = 0.008597140044579702

Write-Host "Applied 279 cell updates"
